# Estadisticos Segundo Parcial 26 Mayo
$wb = $excel.ActiveWorkbook

# --- Calificaciones: enter 2P grades for "Conciencia historica I." (column K)
# and refresh a handful of Final-column (AA) averages for the same subject.
$cal = $wb.Worksheets.Item("Calificaciones")

$k2p = @{
    4 = 10; 5 = 10; 6 = 10; 7 = 6; 8 = 9; 9 = 10; 10 = 9; 11 = 10; 12 = 10;
    13 = 10; 14 = 9; 15 = 5; 16 = 10; 17 = 10; 18 = 5; 19 = 7; 20 = 9; 21 = 7;
    22 = 10; 23 = 10; 24 = 7; 25 = 10; 26 = 5; 27 = 10; 28 = 10; 29 = 6; 30 = 9;
    31 = 10; 32 = 10; 33 = 10; 34 = 6; 35 = 7; 36 = 9; 37 = 5
}
foreach ($row in $k2p.Keys) {
    $cal.Cells.Item($row, 11).Value = $k2p[$row]
}

$aaFinal = @{
    6 = 9; 9 = 8; 10 = 9; 11 = 9; 12 = 9; 13 = 9; 15 = 7; 16 = 9; 17 = 9;
    18 = 8; 20 = 7; 22 = 10; 23 = 8; 24 = 7; 25 = 10; 26 = 5; 28 = 9; 30 = 9;
    32 = 9; 36 = 8; 37 = 7
}
foreach ($row in $aaFinal.Keys) {
    $cal.Cells.Item($row, 27).Value = $aaFinal[$row]
}

# --- Totales: updated average for "Conciencia historica I."
$tot = $wb.Worksheets.Item("Totales")
$tot.Range("H3").Value = 8.4

# --- Rescatables: swap out student "MARIN RODRIGUEZ ABRIL" for
# "PEREZ PEREZ ARELI DANAE" on row 7.
$resc = $wb.Worksheets.Item("Rescatables")
$resc.Range("A7").Value = 23330051920103
$resc.Range("B7").Value = "PEREZ"
$resc.Range("C7").Value = "PEREZ"
$resc.Range("D7").Value = "ARELI DANAE"
